# Update cryptocurrency price/volume figures (scheduled GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text so the stored value matches exactly (e.g. "575.10", not 575.1).
$textCells = [ordered]@{
    'D5' = '575.10'
    'D6' = '171.21'
    'D7' = '0.999'
    'D20' = '16.53'
    'D21' = '486.76'
    'D22' = '7.68'
    'D24' = '82.41'
    'D25' = '12.65'
    'D27' = '10.12'
    'D29' = '7.88'
    'D30' = '2.25'
    'D32' = '27.75'
    'D36' = '48.20'
    'D46' = '134.68'
    'D47' = '365.19'
    'D49' = '24.31'
    'D50' = '2.16'
}
foreach ($cell in $textCells.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $textCells[$cell]
    $rng.Style = "Normal"
}

# Remaining cells are safe to assign directly (not interpreted as numbers).
$plainCells = [ordered]@{
    'D2' = '66.420.12'
    'E2' = '  -0.66%  '
    'D3' = '3.077.83'
    'E3' = '  -1.32%  '
    'E4' = '  +0.01%  '
    'E5' = '  -0.47%  '
    'E6' = '  -0.65%  '
    'E7' = '  -0.08%  '
    'D8' = '3.075.55'
    'E9' = '  -2.07%  '
    'E10' = '  -1.65%  '
    'E11' = '  -2.15%  '
    'E12' = '  -2.52%  '
    'E13' = '  -3.78%  '
    'E14' = '  -3.67%  '
    'D16' = '3.589.62'
    'E16' = '  -1.27%  '
    'D17' = '66.356.43'
    'E17' = '  -0.74%  '
    'E18' = '  -2.67%  '
    'D19' = '3.077.12'
    'E19' = '  -1.34%  '
    'E20' = '  +1.77%  '
    'E21' = '  +2.62%  '
    'E22' = '  -2.47%  '
    'E23' = '  -3.23%  '
    'E24' = '  -1.53%  '
    'E25' = '  -4.53%  '
    'E26' = '  -2.71%  '
    'E27' = '  -1.29%  '
    'E28' = '  +0.02%  '
    'E29' = '  -0.43%  '
    'E30' = '  -4.87%  '
    'E31' = '  -3.07%  '
    'E32' = '  -2.81%  '
    'E33' = '  -3.62%  '
    'E34' = '  -4.21%  '
    'E35' = '  -0.01%  '
    'E36' = '  +2.48%  '
    'E37' = '  -4.65%  '
    'E38' = '  -3.36%  '
    'E40' = '  -3.20%  '
    'E41' = '  -4.67%  '
    'E42' = '  -4.36%  '
    'D43' = '2.770.41'
    'E43' = '  -1.56%  '
    'E44' = '  -1.10%  '
    'E45' = '  -2.79%  '
    'E47' = '  -4.52%  '
    'E48' = '  +0.00%  '
    'E49' = '  -2.49%  '
    'E50' = '  -1.93%  '
    'E51' = '  -2.07%  '
}
foreach ($cell in $plainCells.Keys) {
    $ws.Range($cell).Value = $plainCells[$cell]
}

